$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Ano" header column with the same style as the other headers (A1:C1)
$ws.Range("D1").Value = "Ano"
$ws.Range("D1").Style = $ws.Range("A1").Style

# Fill D2:D10 with the reference year range
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = "2023/2012"
}
